# Updated symbol list on Mon Dec 26 11:19:12 UTC 2022 with GitHub Actions
#
# Refresh of the crypto price/listing table on Sheet1. Most rows just get a
# refreshed Price (column D). Rows 13-26 shifted up one rank (each row took
# on the Coin/Link/Price/Volume info of the row that used to be below it),
# with ProBitToken's refreshed data landing back in at the bottom (row 26).
# A couple of rows also had their Volume(1h) rank-label (column E) text
# tweaked (e.g. "Bestin24h"/"Worstin24h" suffix moved to a different row).
#
# Numeric-looking values are written with a leading apostrophe so Excel
# keeps them as literal text (matching the sheet's existing convention of
# storing Price/Volume as text strings instead of numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}

# --- Simple price-only refreshes (column D) ---
Set-TextCell 2  4 "243.42"
Set-TextCell 4  4 "5.398"
Set-TextCell 5  4 "0.05915"
Set-TextCell 6  4 "3.456"
Set-TextCell 7  4 "6.560"
Set-TextCell 8  4 "0.8122"
Set-TextCell 9  4 "0.9121"
Set-TextCell 10 4 "0.1413"
Set-TextCell 11 4 "0.07426"
Set-TextCell 12 4 "0.03262"

# --- Rows 13-26: rank shift (each row adopts the next row's Coin/Link/Price/
#     Volume label), with refreshed ProBitToken data settling at row 26 ---
Set-TextCell 13 2 "BitrueCoin"
Set-TextCell 13 3 "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell 13 4 "0.03068"
Set-TextCell 13 5 "12BitrueCoinBTR"

Set-TextCell 14 2 "BitMartToken"
Set-TextCell 14 3 "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell 14 4 "0.09339"
Set-TextCell 14 5 "13BitMartTokenBMX"

Set-TextCell 15 2 "MCDex"
Set-TextCell 15 3 "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell 15 4 "3.853"
Set-TextCell 15 5 "14MCDexMCB"

Set-TextCell 16 2 "BitForexToken"
Set-TextCell 16 3 "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell 16 4 "0.001558"
Set-TextCell 16 5 "15BitForexTokenBF"

Set-TextCell 17 2 "CoinExToken"
Set-TextCell 17 3 "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell 17 4 "0.04688"
Set-TextCell 17 5 "16CoinExTokenCET"

Set-TextCell 18 2 "One"
Set-TextCell 18 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell 18 4 "0.01129"
Set-TextCell 18 5 "17OneONEBestin24h"

Set-TextCell 19 2 "TigerCash"
Set-TextCell 19 3 "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell 19 4 "0.005916"
Set-TextCell 19 5 "18TigerCashTCH"

Set-TextCell 20 2 "HotbitToken"
Set-TextCell 20 3 "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell 20 4 "0.004982"
Set-TextCell 20 5 "19HotbitTokenHTB"

Set-TextCell 21 2 "BitKan"
Set-TextCell 21 3 "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell 21 4 "0.0009839"
Set-TextCell 21 5 "20BitKanKAN"

Set-TextCell 22 2 "NitroEx"
Set-TextCell 22 3 "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell 22 4 "0.00008603"
Set-TextCell 22 5 "21NitroExNTX"

Set-TextCell 23 2 "LEO"
Set-TextCell 23 3 "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell 23 4 "3.605"
Set-TextCell 23 5 "22LEOLEO"

Set-TextCell 24 2 "BTSEToken"
Set-TextCell 24 3 "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell 24 4 "2.151"
Set-TextCell 24 5 "23BTSETokenBTSE"

Set-TextCell 25 2 "BitpandaEcosystemToken"
Set-TextCell 25 3 "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell 25 4 "0.3241"
Set-TextCell 25 5 "24BitpandaEcosystemTokenBEST"

Set-TextCell 26 2 "ProBitToken"
Set-TextCell 26 3 "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell 26 4 "0.1324"
Set-TextCell 26 5 "25ProBitTokenPROB"

# --- More simple price-only refreshes (column D) ---
Set-TextCell 27 4 "0.0002901"
Set-TextCell 40 4 "0.03959"
Set-TextCell 41 4 "0.006192"
Set-TextCell 42 4 "0.1073"

# Row 43: price refresh + rank-label text tweak (drop "Worstin24h" suffix)
Set-TextCell 43 4 "0.003001"
Set-TextCell 43 5 "42CEJICEJI"

Set-TextCell 44 4 "0.008926"

# Row 47: price refresh + rank-label text tweak (add "Worstin24h" suffix)
Set-TextCell 47 4 "0.8093"
Set-TextCell 47 5 "46CoinbaseStockTokenCOINWorstin24h"

Set-TextCell 48 4 "0.002330"
